$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("D7").Value = 0.04691346051182279
$ws.Range("E7").Value = 0.1865129431110523
$ws.Range("F7").Value = 0.2643682770640485
$ws.Range("G7").Value = 0.2165951534818422
$ws.Range("H7").Value = 31.91825104663496
$ws.Range("L7").Value = 0.02885241618130121
$ws.Range("M7").Value = 0.1334341014862883
$ws.Range("N7").Value = 0.1803507468632586
$ws.Range("O7").Value = 0.1698599899367159
$ws.Range("P7").Value = 20.41689859926331

$ws = $wb.Worksheets.Item(2)
$ws.Range("D7").Value = 0.03287235298010502
$ws.Range("E7").Value = 0.1193170587321527
$ws.Range("F7").Value = 0.1710301538722946
$ws.Range("G7").Value = 0.1813073439773056
$ws.Range("H7").Value = 18.84215654103123
$ws.Range("L7").Value = 0.03241142083640988
$ws.Range("M7").Value = 0.1184776699774059
$ws.Range("N7").Value = 0.1583014124449199
$ws.Range("O7").Value = 0.1800317217503901
$ws.Range("P7").Value = 17.60310688830728

$ws = $wb.Worksheets.Item(3)
$ws.Range("D7").Value = 0.03318441325406175
$ws.Range("E7").Value = 0.1252566644950564
$ws.Range("F7").Value = 0.1683210465241228
$ws.Range("G7").Value = 0.1821658948707517
$ws.Range("H7").Value = 18.32138605344866
$ws.Range("L7").Value = 0.03796050833193022
$ws.Range("M7").Value = 0.1360457340208823
$ws.Range("N7").Value = 0.1916771712754864
$ws.Range("O7").Value = 0.1948345665736196
$ws.Range("P7").Value = 20.06547374668843

$ws = $wb.Worksheets.Item(4)
$ws.Range("D7").Value = 0.1598974865623922
$ws.Range("E7").Value = 0.3031012090569071
$ws.Range("F7").Value = 0.105392530394474
$ws.Range("G7").Value = 0.3998718376710121
$ws.Range("H7").Value = 11.41564632551584
$ws.Range("L7").Value = 0.3044080361402858
$ws.Range("M7").Value = 0.4577233799393349
$ws.Range("N7").Value = 0.1619838810195003
$ws.Range("O7").Value = 0.5517318516637278
$ws.Range("P7").Value = 18.28631078947365

$ws = $wb.Worksheets.Item(5)
$ws.Range("D7").Value = 0.2598561042740906
$ws.Range("E7").Value = 0.393082015214089
$ws.Range("F7").Value = 0.1672503804963263
$ws.Range("G7").Value = 0.5097608304627678
$ws.Range("H7").Value = 15.93818039721589
$ws.Range("L7").Value = 0.1486769443332642
$ws.Range("M7").Value = 0.2709992587839502
$ws.Range("N7").Value = 0.1089293902373044
$ws.Range("O7").Value = 0.3855864939715397
$ws.Range("P7").Value = 11.06225236777087

$ws = $wb.Worksheets.Item(6)
$ws.Range("D7").Value = 0.1934730861495593
$ws.Range("E7").Value = 0.3508298263297101
$ws.Range("F7").Value = 0.1257689352277472
$ws.Range("G7").Value = 0.4398557560718733
$ws.Range("H7").Value = 13.53529625332182
$ws.Range("L7").Value = 0.1522180623770018
$ws.Range("M7").Value = 0.286408596400033
$ws.Range("N7").Value = 0.1002288701696343
$ws.Range("O7").Value = 0.3901513326608047
$ws.Range("P7").Value = 10.95373872239175

$ws = $wb.Worksheets.Item(7)
$ws.Range("D7").Value = 0.5260366833872417
$ws.Range("E7").Value = 0.5705673325065613
$ws.Range("F7").Value = 1827365085516500
$ws.Range("G7").Value = 0.72528386400584
$ws.Range("H7").Value = 101.6327086128056
$ws.Range("L7").Value = 0.5170865583208656
$ws.Range("M7").Value = 0.6189909348566051
$ws.Range("N7").Value = 1189458637335924
$ws.Range("O7").Value = 0.7190873092475388
$ws.Range("P7").Value = 140.0657036784634

$ws = $wb.Worksheets.Item(8)
$ws.Range("D7").Value = 1.267042639819518
$ws.Range("E7").Value = 0.934943440663304
$ws.Range("F7").Value = 1057520239090499
$ws.Range("G7").Value = 1.125629885805951
$ws.Range("H7").Value = 138.4187521429556
$ws.Range("L7").Value = 1.396843910230323
$ws.Range("M7").Value = 0.9445377265692172
$ws.Range("N7").Value = 917663742709628.8
$ws.Range("O7").Value = 1.181881512771193
$ws.Range("P7").Value = 128.1432478721961

$ws = $wb.Worksheets.Item(9)
$ws.Range("D7").Value = 0.3037978710989326
$ws.Range("E7").Value = 0.461187361244236
$ws.Range("F7").Value = 995708078678606.6
$ws.Range("G7").Value = 0.5511786199581153
$ws.Range("H7").Value = 97.02342091475758
$ws.Range("L7").Value = 0.3104558838041753
$ws.Range("M7").Value = 0.4565850693909362
$ws.Range("N7").Value = 778604034187215.8
$ws.Range("O7").Value = 0.5571856816216434
$ws.Range("P7").Value = 107.2049613212096

$ws = $wb.Worksheets.Item(10)
$ws.Range("D7").Value = 2.07200103491992
$ws.Range("E7").Value = 1.187851047728577
$ws.Range("F7").Value = 1248114887416049
$ws.Range("G7").Value = 1.439444696721593
$ws.Range("H7").Value = 126.6634220714332
$ws.Range("L7").Value = 3.334926492421043
$ws.Range("M7").Value = 1.532943886976013
$ws.Range("N7").Value = 1333255582114454
$ws.Range("O7").Value = 1.826178110815329
$ws.Range("P7").Value = 144.2668875280528

$ws = $wb.Worksheets.Item(11)
$ws.Range("D7").Value = 3.224930151253239
$ws.Range("E7").Value = 1.56809450850699
$ws.Range("F7").Value = 2805591122406188
$ws.Range("G7").Value = 1.795809052002255
$ws.Range("H7").Value = 142.0068760354083
$ws.Range("L7").Value = 3.415090021873789
$ws.Range("M7").Value = 1.391320549173342
$ws.Range("N7").Value = 2687705635541630
$ws.Range("O7").Value = 1.847996218035575
$ws.Range("P7").Value = 121.0415047878117

$ws = $wb.Worksheets.Item(12)
$ws.Range("D7").Value = 0.5743855400463201
$ws.Range("E7").Value = 0.6393050843728957
$ws.Range("F7").Value = 1155788721451477
$ws.Range("G7").Value = 0.7578822732102395
$ws.Range("H7").Value = 148.4429963237464
$ws.Range("L7").Value = 1.074940472317289
$ws.Range("M7").Value = 0.9264864060084336
$ws.Range("N7").Value = 1767106575282857
$ws.Range("O7").Value = 1.03679336047126
$ws.Range("P7").Value = 156.3590355411109

$ws = $wb.Worksheets.Item(13)
$ws.Range("D7").Value = 0.0625997910994754
$ws.Range("E7").Value = 0.1531860723716411
$ws.Range("F7").Value = 1.387247866103559
$ws.Range("G7").Value = 0.2501995025963789
$ws.Range("H7").Value = 29.49183187031828
$ws.Range("L7").Value = 0.06106904735098521
$ws.Range("M7").Value = 0.156049111848684
$ws.Range("N7").Value = 1.365485768840762
$ws.Range("O7").Value = 0.2471215234474432
$ws.Range("P7").Value = 30.12491575558443

$ws = $wb.Worksheets.Item(14)
$ws.Range("D7").Value = 0.05118941941207773
$ws.Range("E7").Value = 0.1273370174838717
$ws.Range("F7").Value = 1.321598690600179
$ws.Range("G7").Value = 0.226250788754598
$ws.Range("H7").Value = 25.40619142710907
$ws.Range("L7").Value = 0.05123666249868912
$ws.Range("M7").Value = 0.1288673666514994
$ws.Range("N7").Value = 1.320984477050704
$ws.Range("O7").Value = 0.2263551689241691
$ws.Range("P7").Value = 25.61218087449125

$ws = $wb.Worksheets.Item(15)
$ws.Range("D7").Value = 0.04000674610406072
$ws.Range("E7").Value = 0.1158203960128188
$ws.Range("F7").Value = 1.22911721571609
$ws.Range("G7").Value = 0.2000168645491193
$ws.Range("H7").Value = 23.842678626271
$ws.Range("L7").Value = 0.03907460257366401
$ws.Range("M7").Value = 0.1136154213718629
$ws.Range("N7").Value = 1.218614462923644
$ws.Range("O7").Value = 0.197672968748041
$ws.Range("P7").Value = 23.39720696724434

$ws = $wb.Worksheets.Item(16)
$ws.Range("D7").Value = 0.5571122454997859
$ws.Range("E7").Value = 0.3290837317365065
$ws.Range("F7").Value = 0.5893408455031428
$ws.Range("G7").Value = 0.7463995213689422
$ws.Range("H7").Value = 16.50080709725727
$ws.Range("L7").Value = 0.5836818109352441
$ws.Range("M7").Value = 0.312714877982272
$ws.Range("N7").Value = 0.6026950289361571
$ws.Range("O7").Value = 0.7639907139064218
$ws.Range("P7").Value = 15.77144523923839

$ws = $wb.Worksheets.Item(17)
$ws.Range("D7").Value = 0.5155407986192199
$ws.Range("E7").Value = 0.3006019845877092
$ws.Range("F7").Value = 0.5959316637386596
$ws.Range("G7").Value = 0.7180116981075029
$ws.Range("H7").Value = 15.26681300847412
$ws.Range("L7").Value = 0.57239993670042
$ws.Range("M7").Value = 0.3089059745282263
$ws.Range("N7").Value = 0.623098472661228
$ws.Range("O7").Value = 0.7565711709419147
$ws.Range("P7").Value = 15.64152761325156

$ws = $wb.Worksheets.Item(18)
$ws.Range("D7").Value = 0.3874259376212878
$ws.Range("E7").Value = 0.2777114738634539
$ws.Range("F7").Value = 0.2593953539199081
$ws.Range("G7").Value = 0.6224354887225565
$ws.Range("H7").Value = 14.13148955516594
$ws.Range("L7").Value = 0.4287089854629551
$ws.Range("M7").Value = 0.2854500612082672
$ws.Range("N7").Value = 0.273780828100897
$ws.Range("O7").Value = 0.6547587230903877
$ws.Range("P7").Value = 14.46699999592231

$ws = $wb.Worksheets.Item(19)
$ws.Range("D7").Value = 0.1386605782654952
$ws.Range("E7").Value = 0.3124623015815102
$ws.Range("F7").Value = 0.4461094991200831
$ws.Range("G7").Value = 0.3723715594208226
$ws.Range("H7").Value = 45.21348982288627

$ws = $wb.Worksheets.Item(20)
$ws.Range("D7").Value = 0.2467826283452084
$ws.Range("E7").Value = 0.388514580633978
$ws.Range("F7").Value = 0.6149863109339269
$ws.Range("G7").Value = 0.496772209715085
$ws.Range("H7").Value = 74.62512569230762

$ws = $wb.Worksheets.Item(21)
$ws.Range("D7").Value = 1.524633521049598
$ws.Range("E7").Value = 0.9061769672805832
$ws.Range("F7").Value = 1.322648818574927
$ws.Range("G7").Value = 1.234760511617373
$ws.Range("H7").Value = 108.378041188903

$ws = $wb.Worksheets.Item(22)
$ws.Range("D7").Value = 1.086878532347106
$ws.Range("E7").Value = 0.7143386181492126
$ws.Range("F7").Value = 0.2690915174269095
$ws.Range("G7").Value = 1.042534667215966
$ws.Range("H7").Value = 30.17736775151323

$ws = $wb.Worksheets.Item(23)
$ws.Range("D7").Value = 13.47039894701386
$ws.Range("E7").Value = 3.187357589496772
$ws.Range("F7").Value = 1.210081305327064
$ws.Range("G7").Value = 3.670204210532959
$ws.Range("H7").Value = 124.4510893863022

$ws = $wb.Worksheets.Item(24)
$ws.Range("D7").Value = 20.586938048087
$ws.Range("E7").Value = 3.582806496431216
$ws.Range("F7").Value = 1.408202573752562
$ws.Range("G7").Value = 4.537283113063037
$ws.Range("H7").Value = 101.8058318287786
